$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 128, shifting existing rows 128-149 down to 130-151
$ws.Range("A128:T129").EntireRow.Insert()

# Fill in the two new rows (128, 129) with a new weekly price report
$ws.Range("A128").Value = 4
$ws.Range("B128").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C128").Value = "Los Lagos"
$ws.Range("D128").Value = 44505
$ws.Range("D128").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E128").Value = 10
$ws.Range("F128").Value = "Fruta"
$ws.Range("G128").Value = 100101
$ws.Range("H128").Value = "Berries"
$ws.Range("I128").Value = 100101007
$ws.Range("J128").Value = "Kiwi"
$ws.Range("K128").Value = "Hayward"
$ws.Range("L128").Value = "Especial"
$ws.Range("M128").Value = 300
$ws.Range("N128").Value = 23000
$ws.Range("O128").Value = 23000
$ws.Range("P128").Value = 23000
$ws.Range("Q128").Value = "$/caja 15 kilos"
$ws.Range("R128").Value = "Provincia de Curicó"
$ws.Range("S128").Value = 1533
$ws.Range("T128").Value = 15

$ws.Range("A129").Value = 4
$ws.Range("B129").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C129").Value = "Los Lagos"
$ws.Range("D129").Value = 44505
$ws.Range("D129").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E129").Value = 10
$ws.Range("F129").Value = "Fruta"
$ws.Range("G129").Value = 100101
$ws.Range("H129").Value = "Berries"
$ws.Range("I129").Value = 100101007
$ws.Range("J129").Value = "Kiwi"
$ws.Range("K129").Value = "Hayward"
$ws.Range("L129").Value = "Primera"
$ws.Range("M129").Value = 600
$ws.Range("N129").Value = 16000
$ws.Range("O129").Value = 17000
$ws.Range("P129").Value = 16500
$ws.Range("Q129").Value = "$/caja 15 kilos"
$ws.Range("R129").Value = "Provincia de Curicó"
$ws.Range("S129").Value = 1100
$ws.Range("T129").Value = 15
